# Remove all old terms for MDR introduction.
# - Rename "start_mdr_introduce_time" (row 5) to "mdr_introduce_time".
# - Remove the "end_mdr_introduce_time" row (row 6) entirely, shifting the
#   rows below it (susceptible_fully, active, age_breakpoints) up by one.
# - Leave the final selection on B9 (matches the author's saved selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the MDR introduction parameter - it no longer needs the "start_" prefix
# since the corresponding "end_mdr_introduce_time" row is being removed.
$ws.Range("A5").Value = "mdr_introduce_time"

# Delete the whole "end_mdr_introduce_time" row (old row 6). This shifts
# susceptible_fully/active/age_breakpoints up one row and automatically
# adjusts the dependent data-validation ranges.
$ws.Rows("6:6").Delete()

# Match the workbook's saved cursor position after the edit.
$ws.Range("B9").Select() | Out-Null
